$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 44: a new monthly data point (01-07-2021) appended after the
# existing 01-06-2021 row.
#
# Column A stores these month labels as plain text shared strings with
# the sheet's default (unstyled) cell format - same as every row above.
# Assigning a date-shaped literal ("01-07-2021") straight to Range.Value
# makes the COM layer's type inference treat it as a date, converting it
# to a serial number and minting/applying a new date-formatted style -
# which would incorrectly touch styles.xml and store a numeric date
# instead of the original text. Writing it as a text formula first and
# then copy/paste-special-ing the formula's cached text result back onto
# itself converts the cell to a plain literal text value (a shared
# string) without ever invoking the Value-setter's date heuristic, so no
# style is created and styles.xml stays untouched.
$ws.Range("A44").Formula = '="01-07-2021"'
$ws.Range("A44").Copy()
$ws.Range("A44").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B44").Value = 109.76
$ws.Range("C44").Value = 108.01
$ws.Range("D44").Value = 111.4
$ws.Range("E44").Value = 107.81
$ws.Range("F44").Value = 118.06
